$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet1: insert a new "portrait indicator" row right after row 18
#     (SP_BUSINESS_PACKAGE), shifting the existing rows 19..107 down to 20..108.
#     Copy bottom-up so every row keeps its original value + formatting.
for ($r = 107; $r -ge 19; $r--) {
    $src = $ws.Range("A" + $r + ":B" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":B" + ($r + 1))
    $src.Copy($dst)
}

# New row 19 content
$ws.Cells.Item(19, 1).Value = "SP_BIG_PLANS"
$ws.Cells.Item(19, 2).Value = "Доля клиентов с ПУ Большие планы"

# Match the look of the other "SP_*" portrait-indicator rows: plain style in
# column A (same as its neighbours), highlighted style in column B (same as
# e.g. row 18's description cell).
$ws.Cells.Item(20, 1).Copy()
$ws.Cells.Item(19, 1).PasteSpecial(-4122)
$ws.Cells.Item(18, 2).Copy()
$ws.Cells.Item(19, 2).PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Selection / view: jump back to the top and land on the new cell
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("B19").Select()
